$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 357-477 (1-indexed) of the "Plátano" dataset: the logged prices for
# the newest period were prepended as rows 357-358, pushing the rest of the
# series down by two rows; the two rows that fell off the bottom of the
# previous range are appended as new rows 476-477.
$rows = @(
  @{R=357; D=44559; L='Pintón'; M=1580; N=10000; O=10000; P=10000; S=500},
  @{R=358; D=44559; L='Primera Pintón'; M=560; N=11000; O=11000; P=11000; S=550},
  @{R=359; D=44341; L='Pintón'; M=900; N=10000; O=10000; P=10000; S=500},
  @{R=360; D=44341; L='Primera Pintón'; M=550; N=11000; O=11000; P=11000; S=550},
  @{R=361; D=44286; L='Pintón'; M=800; N=13000; O=13000; P=13000; S=650},
  @{R=362; D=44279; L='Pintón'; M=600; N=12000; O=12000; P=12000; S=600},
  @{R=363; D=44279; L='Primera Pintón'; M=230; N=13000; O=13000; P=13000; S=650},
  @{R=364; D=44208; L='Pintón'; M=800; N=12000; O=12000; P=12000; S=600},
  @{R=365; D=44208; L='Primera Pintón'; M=250; N=13000; O=13000; P=13000; S=650},
  @{R=366; D=44264; L='Pintón'; M=800; N=14000; O=14000; P=14000; S=700},
  @{R=367; D=44264; L='Primera Pintón'; M=300; N=15000; O=15000; P=15000; S=750},
  @{R=368; D=44322; L='Pintón'; M=700; N=15000; O=15000; P=15000; S=750},
  @{R=369; D=44322; L='Primera Pintón'; M=300; N=16000; O=16000; P=16000; S=800},
  @{R=370; D=44491; L='Pintón'; M=1140; N=22000; O=23000; P=22526; S=1126},
  @{R=371; D=44389; L='Pintón'; M=300; N=10000; O=10000; P=10000; S=500},
  @{R=372; D=44389; L='Primera Pintón'; M=450; N=12000; O=12000; P=12000; S=600},
  @{R=373; D=44391; L='Pintón'; M=300; N=9000; O=9000; P=9000; S=450},
  @{R=374; D=44391; L='Primera Pintón'; M=500; N=11000; O=11000; P=11000; S=550},
  @{R=375; D=44396; L='Pintón'; M=800; N=11000; O=11000; P=11000; S=550},
  @{R=376; D=44396; L='Primera Pintón'; M=400; N=12000; O=12000; P=12000; S=600},
  @{R=377; D=44510; L='Pintón'; M=800; N=15000; O=15000; P=15000; S=750},
  @{R=378; D=44510; L='Primera Pintón'; M=500; N=17000; O=17000; P=17000; S=850},
  @{R=379; D=44232; L='Pintón'; M=840; N=13000; O=13000; P=13000; S=650},
  @{R=380; D=44232; L='Primera Pintón'; M=310; N=16000; O=16000; P=16000; S=800},
  @{R=381; D=44551; L='Pintón'; M=900; N=10000; O=10000; P=10000; S=500},
  @{R=382; D=44551; L='Primera Pintón'; M=1000; N=11000; O=12000; P=11800; S=590},
  @{R=383; D=44386; L='Pintón'; M=900; N=10000; O=10000; P=10000; S=500},
  @{R=384; D=44386; L='Primera Pintón'; M=300; N=11000; O=11000; P=11000; S=550},
  @{R=385; D=44519; L='Pintón'; M=520; N=14000; O=14000; P=14000; S=700},
  @{R=386; D=44519; L='Primera Pintón'; M=300; N=16000; O=16000; P=16000; S=800},
  @{R=387; D=44420; L='Pintón'; M=520; N=10000; O=10000; P=10000; S=500},
  @{R=388; D=44420; L='Primera Pintón'; M=300; N=11000; O=11000; P=11000; S=550},
  @{R=389; D=44414; L='Pintón'; M=500; N=11000; O=11000; P=11000; S=550},
  @{R=390; D=44414; L='Primera Pintón'; M=300; N=12000; O=12000; P=12000; S=600},
  @{R=391; D=44543; L='Pintón'; M=1050; N=12000; O=12000; P=12000; S=600},
  @{R=392; D=44543; L='Primera Pintón'; M=480; N=13000; O=13000; P=13000; S=650},
  @{R=393; D=44321; L='Pintón'; M=850; N=15000; O=15000; P=15000; S=750},
  @{R=394; D=44321; L='Primera Pintón'; M=350; N=17000; O=17000; P=17000; S=850},
  @{R=395; D=44385; L='Pintón'; M=800; N=9000; O=9000; P=9000; S=450},
  @{R=396; D=44385; L='Primera Pintón'; M=500; N=10000; O=10000; P=10000; S=500},
  @{R=397; D=44278; L='Pintón'; M=900; N=11000; O=11000; P=11000; S=550},
  @{R=398; D=44278; L='Primera Pintón'; M=150; N=12000; O=12000; P=12000; S=600},
  @{R=399; D=44308; L='Pintón'; M=900; N=10500; O=10500; P=10500; S=525},
  @{R=400; D=44308; L='Primera Pintón'; M=540; N=12000; O=12000; P=12000; S=600},
  @{R=401; D=44281; L='Pintón'; M=300; N=11000; O=11000; P=11000; S=550},
  @{R=402; D=44281; L='Primera Pintón'; M=180; N=13000; O=13000; P=13000; S=650},
  @{R=403; D=44187; L='Pintón'; M=800; N=10000; O=10000; P=10000; S=500},
  @{R=404; D=44187; L='Primera Pintón'; M=300; N=11000; O=11000; P=11000; S=550},
  @{R=405; D=44474; L='Pintón'; M=800; N=18000; O=18000; P=18000; S=900},
  @{R=406; D=44474; L='Primera Pintón'; M=500; N=19000; O=19000; P=19000; S=950},
  @{R=407; D=44446; L='Pintón'; M=800; N=18000; O=18000; P=18000; S=900},
  @{R=408; D=44446; L='Primera Pintón'; M=600; N=19000; O=19000; P=19000; S=950},
  @{R=409; D=44350; L='Primera Pintón'; M=260; N=11000; O=11000; P=11000; S=550},
  @{R=410; D=44529; L='Pintón'; M=1000; N=18000; O=18000; P=18000; S=900},
  @{R=411; D=44529; L='Primera'; M=600; N=20000; O=20000; P=20000; S=1000},
  @{R=412; D=44405; L='Pintón'; M=800; N=14000; O=14000; P=14000; S=700},
  @{R=413; D=44405; L='Primera Pintón'; M=450; N=16000; O=16000; P=16000; S=800},
  @{R=414; D=44413; L='Pintón'; M=540; N=11000; O=11000; P=11000; S=550},
  @{R=415; D=44413; L='Primera Pintón'; M=540; N=12000; O=12000; P=12000; S=600},
  @{R=416; D=44238; L='Pintón'; M=400; N=7000; O=7000; P=7000; S=350},
  @{R=417; D=44238; L='Primera Pintón'; M=600; N=8000; O=8000; P=8000; S=400},
  @{R=418; D=44257; L='Pintón'; M=800; N=14000; O=14000; P=14000; S=700},
  @{R=419; D=44411; L='Pintón'; M=1050; N=12000; O=12000; P=12000; S=600},
  @{R=420; D=44411; L='Primera Pintón'; M=840; N=13000; O=14000; P=13357; S=668},
  @{R=421; D=44175; L='Pintón'; M=800; N=13000; O=13000; P=13000; S=650},
  @{R=422; D=44175; L='Primera Pintón'; M=500; N=14000; O=14000; P=14000; S=700},
  @{R=423; D=44196; L='Primera'; M=600; N=12000; O=12000; P=12000; S=600},
  @{R=424; D=44196; L='Primera Pintón'; M=320; N=14000; O=14000; P=14000; S=700},
  @{R=425; D=44200; L='Pintón'; M=600; N=14000; O=14000; P=14000; S=700},
  @{R=426; D=44459; L='Pintón'; M=700; N=17000; O=17000; P=17000; S=850},
  @{R=427; D=44459; L='Primera Pintón'; M=400; N=18000; O=18000; P=18000; S=900},
  @{R=428; D=44188; L='Pintón'; M=550; N=10000; O=10000; P=10000; S=500},
  @{R=429; D=44188; L='Primera Pintón'; M=300; N=11000; O=11000; P=11000; S=550},
  @{R=430; D=44258; L='Pintón'; M=850; N=14000; O=14000; P=14000; S=700},
  @{R=431; D=44298; L='Pintón'; M=1050; N=11000; O=11000; P=11000; S=550},
  @{R=432; D=44298; L='Primera Pintón'; M=800; N=13000; O=13000; P=13000; S=650},
  @{R=433; D=44432; L='Pintón'; M=500; N=12000; O=12000; P=12000; S=600},
  @{R=434; D=44432; L='Primera Pintón'; M=500; N=13000; O=13000; P=13000; S=650},
  @{R=435; D=44428; L='Pintón'; M=350; N=13000; O=13000; P=13000; S=650},
  @{R=436; D=44428; L='Primera Pintón'; M=800; N=14000; O=14000; P=14000; S=700},
  @{R=437; D=44340; L='Pintón'; M=1250; N=10000; O=10000; P=10000; S=500},
  @{R=438; D=44340; L='Primera Pintón'; M=480; N=11500; O=11500; P=11500; S=575},
  @{R=439; D=44224; L='Pintón'; M=600; N=13000; O=13000; P=13000; S=650},
  @{R=440; D=44224; L='Primera Pintón'; M=400; N=14000; O=14000; P=14000; S=700},
  @{R=441; D=44329; L='Pintón'; M=1000; N=11000; O=11000; P=11000; S=550},
  @{R=442; D=44329; L='Primera Pintón'; M=800; N=12000; O=12000; P=12000; S=600},
  @{R=443; D=44452; L='Pintón'; M=300; N=21000; O=21000; P=21000; S=1050},
  @{R=444; D=44452; L='Primera Pintón'; M=400; N=23000; O=23000; P=23000; S=1150},
  @{R=445; D=44270; L='Pintón'; M=1080; N=11000; O=11000; P=11000; S=550},
  @{R=446; D=44270; L='Primera Pintón'; M=840; N=13000; O=13000; P=13000; S=650},
  @{R=447; D=44195; L='Pintón'; M=800; N=12000; O=12000; P=12000; S=600},
  @{R=448; D=44195; L='Primera Pintón'; M=1000; N=13000; O=14000; P=13500; S=675},
  @{R=449; D=44473; L='Pintón'; M=1000; N=19000; O=19000; P=19000; S=950},
  @{R=450; D=44473; L='Primera Pintón'; M=500; N=20000; O=20000; P=20000; S=1000},
  @{R=451; D=44398; L='Pintón'; M=600; N=14000; O=14000; P=14000; S=700},
  @{R=452; D=44398; L='Primera Pintón'; M=250; N=16000; O=16000; P=16000; S=800},
  @{R=453; D=44536; L='Pintón'; M=450; N=15000; O=15000; P=15000; S=750},
  @{R=454; D=44536; L='Primera Pintón'; M=300; N=16000; O=16000; P=16000; S=800},
  @{R=455; D=44302; L='Pintón'; M=600; N=12000; O=12000; P=12000; S=600},
  @{R=456; D=44302; L='Primera Pintón'; M=300; N=13000; O=13000; P=13000; S=650},
  @{R=457; D=44511; L='Pintón'; M=800; N=16000; O=16000; P=16000; S=800},
  @{R=458; D=44511; L='Primera Pintón'; M=450; N=15000; O=15000; P=15000; S=750},
  @{R=459; D=44239; L='Pintón'; M=400; N=7000; O=7000; P=7000; S=350},
  @{R=460; D=44239; L='Primera Pintón'; M=500; N=8000; O=8000; P=8000; S=400},
  @{R=461; D=44461; L='Pintón'; M=300; N=14000; O=14000; P=14000; S=700},
  @{R=462; D=44461; L='Primera Pintón'; M=550; N=15000; O=15000; P=15000; S=750},
  @{R=463; D=44463; L='Pintón'; M=700; N=13000; O=13000; P=13000; S=650},
  @{R=464; D=44463; L='Primera Pintón'; M=800; N=14000; O=15000; P=14375; S=719},
  @{R=465; D=44357; L='Pintón'; M=800; N=10000; O=10000; P=10000; S=500},
  @{R=466; D=44357; L='Primera Pintón'; M=300; N=11000; O=11000; P=11000; S=550},
  @{R=467; D=44371; L='Pintón'; M=900; N=11000; O=11000; P=11000; S=550},
  @{R=468; D=44371; L='Primera Pintón'; M=360; N=12000; O=12000; P=12000; S=600},
  @{R=469; D=44365; L='Pintón'; M=800; N=15000; O=15000; P=15000; S=750},
  @{R=470; D=44365; L='Primera Pintón'; M=450; N=17000; O=17000; P=17000; S=850},
  @{R=471; D=44194; L='Pintón'; M=300; N=13000; O=13000; P=13000; S=650},
  @{R=472; D=44194; L='Primera Pintón'; M=500; N=12000; O=12000; P=12000; S=600},
  @{R=473; D=44313; L='Pintón'; M=1150; N=12000; O=14000; P=12522; S=626},
  @{R=474; D=44540; L='Pintón'; M=800; N=12500; O=13000; P=12812; S=641},
  @{R=475; D=44540; L='Primera Pintón'; M=500; N=14000; O=14000; P=14000; S=700},
  @{R=476; D=44272; L='Pintón'; M=800; N=11000; O=11000; P=11000; S=550},
  @{R=477; D=44272; L='Primera Pintón'; M=260; N=12000; O=12000; P=12000; S=600}
)

$dateFormat = $ws.Range("D357").NumberFormat

foreach ($row in $rows) {
    $r = $row.R
    if ($r -gt 475) {
        # Brand-new row: populate every column.
        $ws.Cells.Item($r, 1).Value = 5
        $ws.Cells.Item($r, 2).Value = "Macroferia Regional de Talca"
        $ws.Cells.Item($r, 3).Value = "Maule"
        $ws.Cells.Item($r, 5).Value = 7
        $ws.Cells.Item($r, 6).Value = "Fruta"
        $ws.Cells.Item($r, 7).Value = 100108
        $ws.Cells.Item($r, 8).Value = "Tropicales y subtropicales"
        $ws.Cells.Item($r, 9).Value = 100108006
        $ws.Cells.Item($r, 10).Value = "Plátano"
        $ws.Cells.Item($r, 11).Value = "Sin especificar"
        $ws.Cells.Item($r, 17).Value = "`$/caja 20 kilos"
        $ws.Cells.Item($r, 18).Value = "Ecuador"
        $ws.Cells.Item($r, 20).Value = 20
    }

    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 4).NumberFormat = $dateFormat
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 19).Value = $row.S
}

Write-Output "Updated $($rows.Count) rows; new dimension should be A1:T477"
